# Update the "Tool" column (D) order on the "query" worksheet to reflect
# the new folder/tool naming order used by the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Analytics rows (2-17): reorder "Tableau, Power BI, Slack, Jira" values.
$ws.Range("D2").Value  = "Jira, Slack, Tableau, Power BI"
$ws.Range("D3").Value  = "Slack, Jira, Tableau, Power BI"
$ws.Range("D4").Value  = "Jira, Slack, Tableau, Power BI"
$ws.Range("D5").Value  = "Slack, Jira, Tableau, Power BI"
$ws.Range("D6").Value  = "Slack, Jira, Tableau, Power BI"
$ws.Range("D7").Value  = "Slack, Jira, Tableau, Power BI"
$ws.Range("D8").Value  = "Jira, Slack, Tableau, Power BI"
$ws.Range("D9").Value  = "Slack, Jira, Tableau, Power BI"
$ws.Range("D10").Value = "Slack, Jira, Tableau, Power BI"
$ws.Range("D11").Value = "Slack, Jira, Tableau, Power BI"
$ws.Range("D12").Value = "Slack, Jira, Tableau, Power BI"
$ws.Range("D13").Value = "Slack, Jira, Tableau, Power BI"
$ws.Range("D14").Value = "Jira, Slack, Tableau, Power BI"
$ws.Range("D15").Value = "Jira, Slack, Tableau, Power BI"
$ws.Range("D16").Value = "Slack, Jira, Tableau, Power BI"
$ws.Range("D17").Value = "Slack, Jira, Tableau, Power BI"

# Engineering rows (18-37): reorder "GitLab, Bitbucket, Docker, IntelliJ IDEA, Jenkins" values.
for ($r = 18; $r -le 37; $r++) {
    $ws.Range("D$r").Value = "Docker, IntelliJ IDEA, Bitbucket, GitLab, Jenkins"
}
